$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vpa_cell_neg_top_hits_w_FC_pval")

$ws.Range("A22").Value = "Galactitol"
$ws.Range("A23").Value = "Glucose 6-Phosphate"
$ws.Range("A24").Value = "Ribose"
$ws.Range("A28").Value = "Glyceraldehyde 3-Phosphate"
$ws.Range("A29").Value = "Glyceraldehyde 3-Phosphate"
$ws.Range("A33").Value = "Glycerol 1-Phosphoserine"
$ws.Range("A34").Value = "Glycerol 1-Phosphoserine"
$ws.Range("A35").Value = "Glycerol-3-Phosphocholine"
$ws.Range("A46").Value = "Methyl-Lysine"
$ws.Range("A49").Value = "N-Acetylglutamine"
$ws.Range("A59").Value = "Sedoheptulose 7-Phosphate"
$ws.Range("A70").Value = "Valproic Acid"
$ws.Range("A71").Value = "Valproic Acid"

$ws.Range("A72").Select()
